# Dumbwaiter teleport and Intro scene
#  - Made the ladder trigger the dumbwaiter floor as a valid teleport target
#  - Added the Intro scene (text and flow only, sound to come)
#
# The Opening/Intro sections were reworked and a new "End" section was
# appended, which reflows almost every row below it (a blank separator row
# now follows the Intro block, mirroring the blank row that already
# separated Main from the header). Rather than trying to replay every
# row-shift individually, clear the whole previously-used table and
# rewrite it cell-by-cell to match the final layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:D40").ClearContents()

# Header row
$ws.Range("A1").Value = 'Section'
$ws.Range("B1").Value = 'Element'
$ws.Range("C1").Value = 'Description'
$ws.Range("D1").Value = 'Status'

# Opening section (UI) - reworked from Title Page/Sitting-Standing/Start
# button into the new Menu Scene / Exit Function / Intro Scene flow
$ws.Range("A2").Value = 'Opening'
$ws.Range("B2").Value = 'UI'
$ws.Range("C2").Value = 'Menu Scene'
$ws.Range("D2").Value = 'Done'
$ws.Range("B3").Value = 'UI'
$ws.Range("C3").Value = 'Exit Function'
$ws.Range("B4").Value = 'UI'
$ws.Range("C4").Value = 'Intro Scene'
$ws.Range("D4").Value = 'Done'

# Intro section (Audio) - new Narration/VO text, ambient music line renamed
$ws.Range("A5").Value = 'Intro'
$ws.Range("B5").Value = 'Audio'
$ws.Range("C5").Value = 'Narration'
$ws.Range("B6").Value = 'Audio'
$ws.Range("C6").Value = 'Bang of door'
$ws.Range("B7").Value = 'Audio'
$ws.Range("C7").Value = 'VO of owner '
$ws.Range("B8").Value = 'Audio'
$ws.Range("C8").Value = 'Ambient music'
$ws.Range("D8").Value = 'Done'

# (row 9 left blank, as a section separator)

# Main section (Enviro / Inventory / Puzzles / Gameplay / Clues)
$ws.Range("A10").Value = 'Main'
$ws.Range("B10").Value = 'Enviro'
$ws.Range("C10").Value = 'Glasses on shelves'
$ws.Range("D10").Value = 'Done'
$ws.Range("C11").Value = 'Other kitchen wares on shelves'
$ws.Range("D11").Value = 'Done'
$ws.Range("C12").Value = 'Bricks for smashing vases'
$ws.Range("D12").Value = 'Done'
$ws.Range("C13").Value = 'Make all drawers work'
$ws.Range("D13").Value = 'Done'
$ws.Range("C14").Value = 'Handle puzzle locked drawer'
$ws.Range("D14").Value = 'Done'
$ws.Range("C15").Value = 'Handle key locked drawer'
$ws.Range("D15").Value = 'Done'
$ws.Range("C16").Value = 'Make taps'
$ws.Range("D16").Value = 'Done'
$ws.Range("C17").Value = 'Place taps'
$ws.Range("D17").Value = 'Done'
$ws.Range("C18").Value = 'Add FRAGILE box to clue vase smashing'
$ws.Range("D18").Value = 'Done'
$ws.Range("C19").Value = 'Door handle on inside of top door (non functional)'

$ws.Range("B20").Value = 'Inventory'
$ws.Range("C20").Value = 'Camera rendering to plane'
$ws.Range("D20").Value = 'Done'
$ws.Range("C21").Value = 'Highlights for inventory objects'
$ws.Range("D21").Value = 'Done'
$ws.Range("C22").Value = 'Spots for inventory objects'
$ws.Range("D22").Value = 'Done'
$ws.Range("C23").Value = 'Code for selecting inventory objects'
$ws.Range("D23").Value = 'Nearly done - return positioning is off'
$ws.Range("C24").Value = 'Code for storing inventory objects'
$ws.Range("D24").Value = 'Done'

$ws.Range("B25").Value = 'Puzzles'
$ws.Range("C25").Value = 'Wine bottle placement in slots'
$ws.Range("D25").Value = 'Done'
$ws.Range("C26").Value = 'Bricks smashing vases'
$ws.Range("D26").Value = 'Done'
$ws.Range("C27").Value = 'Searching sack for key'
$ws.Range("D27").Value = 'Done'
# Ladder -> dumbwaiter teleport target is now finished
$ws.Range("C28").Value = 'Make ladder necessary to teleport into dumbwaiter'
$ws.Range("D28").Value = 'Done'
$ws.Range("C29").Value = 'Make colour changes work for kegs'
$ws.Range("D29").Value = 'Done'
$ws.Range("C30").Value = 'Make tap insertion work for kegs'
$ws.Range("D30").Value = 'Done'

$ws.Range("B31").Value = 'Gameplay'
$ws.Range("C31").Value = 'Make teleport targets/rotations'
$ws.Range("D31").Value = 'Abandoned - not practical with this layout'
$ws.Range("C32").Value = 'Allow teleporting with full hands'
$ws.Range("D32").Value = 'Done'

$ws.Range("B33").Value = 'UI'
$ws.Range("C33").Value = 'Inventory viewing'
$ws.Range("D33").Value = 'Done'

$ws.Range("B34").Value = 'Clues'
$ws.Range("C34").Value = 'Keg colour clue on barrel at bottom of stairs'
$ws.Range("D34").Value = 'Done'
$ws.Range("C35").Value = 'Keg colour clue inside dumbwaiter'
$ws.Range("D35").Value = 'Done'
$ws.Range("C36").Value = 'Keg colour clue in drawer'
$ws.Range("D36").Value = 'Done'

# (row 37 left blank, as a section separator)

# New End section (Audio / Enviro / UI)
$ws.Range("A38").Value = 'End'
$ws.Range("B38").Value = 'Audio'
$ws.Range("C38").Value = 'Win music'
$ws.Range("D38").Value = 'Done'
$ws.Range("B39").Value = 'Enviro'
$ws.Range("C39").Value = 'Animation  of middle cask opening'
$ws.Range("B40").Value = 'UI'
$ws.Range("C40").Value = 'Credits scene'

# Restore the active-cell selection recorded for the sheet view.
$ws.Range("D6").Select()
